# Update the cryptos symbol list with freshly scraped prices / volumes.
# Values in columns D (Price) and E (Volume(1h)) are stored as plain text
# (not numbers/percentages), and a few rows shifted up by one position
# (GateToken moved from row 19 to row 7, etc.) with their rank-links/prices
# following along. We write every changed cell explicitly below.
#
# NOTE: these columns hold numeric-looking / percentage-looking strings
# that Excel would otherwise auto-convert to numbers (losing formatting
# like trailing zeros, e.g. "44.30" -> 44.3). To keep them as literal text
# we briefly force the cell to Text format before assigning, then clear
# the format again so no stray styling is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $newValue)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value2 = $newValue
    $c.ClearFormats()
}

Set-TextValue 'D2' '326.52'
Set-TextValue 'E2' '0.04%'
Set-TextValue 'D3' '44.30'
Set-TextValue 'E3' '-1.44%'
Set-TextValue 'D4' '5.509'
Set-TextValue 'E4' '-1.16%'
Set-TextValue 'D5' '0.08026'
Set-TextValue 'E5' '-0.76%'
Set-TextValue 'D6' '2.006'
Set-TextValue 'E6' '5.03%'
Set-TextValue 'B7' 'GateToken'
Set-TextValue 'C7' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D7' '4.287'
Set-TextValue 'E7' '-1.13%'
Set-TextValue 'B8' 'BTSEToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue 'D8' '2.567'
Set-TextValue 'E8' '-5.51%'
Set-TextValue 'B9' 'MXToken'
Set-TextValue 'C9' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D9' '0.9486'
Set-TextValue 'E9' '-0.02%'
Set-TextValue 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D10' '0.1148'
Set-TextValue 'E10' '-1.74%'
Set-TextValue 'B11' 'WazirX'
Set-TextValue 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.1847'
Set-TextValue 'E11' '-2.51%'
Set-TextValue 'B12' 'MCDex'
Set-TextValue 'C12' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D12' '12.55'
Set-TextValue 'E12' '47.70%'
Set-TextValue 'B13' 'MandalaExchangeToken'
Set-TextValue 'C13' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D13' '0.09775'
Set-TextValue 'E13' '-3.25%'
Set-TextValue 'B14' 'BitrueCoin'
Set-TextValue 'C14' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D14' '0.04614'
Set-TextValue 'E14' '10.34%'
Set-TextValue 'B15' 'BitMartToken'
Set-TextValue 'C15' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D15' '0.1066'
Set-TextValue 'E15' '0.12%'
Set-TextValue 'B16' 'BitForexToken'
Set-TextValue 'C16' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D16' '0.001274'
Set-TextValue 'E16' '0.26%'
Set-TextValue 'B17' 'CoinExToken'
Set-TextValue 'C17' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue 'D17' '0.04082'
Set-TextValue 'E17' '-4.63%'
Set-TextValue 'B18' 'TigerCash'
Set-TextValue 'C18' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D18' '0.005814'
Set-TextValue 'E18' '-2.86%'
Set-TextValue 'B19' 'LEO'
Set-TextValue 'C19' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D19' '3.367'
Set-TextValue 'E19' '-6.88%'
Set-TextValue 'D20' '0.3480'
Set-TextValue 'E20' '-0.21%'
Set-TextValue 'D21' '0.1406'
Set-TextValue 'E21' '2.32%'
Set-TextValue 'E22' '-4.53%'
Set-TextValue 'D23' '0.001244'
Set-TextValue 'E23' '0.61%'
Set-TextValue 'E24' '-6.14%'
Set-TextValue 'E25' '-3.60%'
Set-TextValue 'E26' '-6.52%'
Set-TextValue 'D38' '0.02558'
Set-TextValue 'E38' '-4.23%'
Set-TextValue 'D39' '0.05572'
Set-TextValue 'E39' '0.30%'
Set-TextValue 'D40' '0.007509'
Set-TextValue 'E40' '-2.83%'
Set-TextValue 'D41' '0.1397'
Set-TextValue 'E41' '0.22%'
Set-TextValue 'D42' '0.007615'
Set-TextValue 'E42' '-32.82%'
Set-TextValue 'D43' '0.002014'
Set-TextValue 'E43' '-2.24%'
Set-TextValue 'D44' '0.008511'
Set-TextValue 'E44' '-2.17%'
Set-TextValue 'D45' '0.00007104'
Set-TextValue 'E45' '-0.35%'
Set-TextValue 'E46' '-0.40%'
Set-TextValue 'E47' '54.88%'
Set-TextValue 'D48' '0.001841'
Set-TextValue 'E48' '-47.42%'
Set-TextValue 'E49' '-0.40%'
Set-TextValue 'E50' '-0.40%'

Write-Output "Applied 94 cell updates to crypto symbol list."
